$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 12525539
$arr[0,1] = 8929384
$arr[0,2] = 20916566
$arr[0,3] = 26788152
$arr[0,4] = 62749698
$arr[0,5] = -26785602
$arr[0,6] = -62754798
$ws.Range("H137:N137").Value = $arr

$ws.Range("H140").Value = 64320
$ws.Range("J140").Value = 64320
$ws.Range("L140").Value = 64320
$ws.Range("N140").Value = -74680

$ws = $wb.Worksheets.Item("ARM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 5143.38
$arr[0,1] = 3298.1797
$arr[0,2] = 20072.727
$arr[0,3] = 3298.1797
$arr[0,4] = 20072.727
$arr[0,5] = -3011.1797
$arr[0,6] = -20646.727
$ws.Range("H32:N32").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1962547.2
$arr[0,1] = 1260.5
$arr[0,2] = 4904477.5
$arr[0,3] = 1260.5
$arr[0,4] = 4904477.5
$arr[0,5] = -1048.5
$arr[0,6] = -4904901.5
$ws.Range("H61:N61").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 37039030
$arr[0,1] = 38462572
$arr[0,2] = 33337830
$arr[0,3] = 38462572
$arr[0,4] = 33337830
$arr[0,5] = -38461698
$arr[0,6] = -33339578
$ws.Range("H74:N74").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 37039030
$arr[0,1] = 38462572
$arr[0,2] = 33337830
$arr[0,3] = 192312860
$arr[0,4] = 166689150
$arr[0,5] = -192308492
$arr[0,6] = -166697886
$ws.Range("H77:N77").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 34487716
$arr[0,1] = 43483736
$arr[0,2] = 2985.1667
$arr[0,3] = 130451208
$arr[0,4] = 8955.500100000001
$arr[0,5] = -130448678
$arr[0,6] = -14015.5001
$ws.Range("H132:N132").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1962547.2
$arr[0,1] = 1260.5
$arr[0,2] = 4904477.5
$arr[0,3] = 3781.5
$arr[0,4] = 14713432.5
$arr[0,5] = -1231.5
$arr[0,6] = -14718532.5
$ws.Range("H136:N136").Value = $arr

$ws.Range("H139").Value = 39630
$ws.Range("J139").Value = 39630
$ws.Range("L139").Value = 39630
$ws.Range("N139").Value = -49910

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1500
$ws.Range("I5").Value = 1500
$ws.Range("K5").Value = 1500
$ws.Range("M5").Value = -1387

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 13904247
$arr[0,1] = 17242432
$arr[0,2] = 74625.71000000001
$arr[0,3] = 51727296
$arr[0,4] = 223877.13
$arr[0,5] = -51724761
$arr[0,6] = -228947.13
$ws.Range("H134:N134").Value = $arr

$ws = $wb.Worksheets.Item("CRP")
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$ws.Range("H7:L7").Value = $arr
$ws.Range("M7:N7").ClearContents()

$ws.Range("H50").Value = 11814.667
$ws.Range("J50").Value = 11814.667
$ws.Range("L50").Value = 11814.667
$ws.Range("N50").Value = -13064.667

$ws.Range("H51").Value = 9172.833000000001
$ws.Range("J51").Value = 9172.833000000001
$ws.Range("L51").Value = 9172.833000000001
$ws.Range("N51").Value = -10644.833

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 970153.5600000001
$arr[0,1] = 4085.3
$arr[0,2] = 2674979.8
$arr[0,3] = 4085.3
$arr[0,4] = 2674979.8
$arr[0,5] = -3882.3
$arr[0,6] = -2675385.8
$ws.Range("H58:N58").Value = $arr

$ws.Range("H59").Value = 15793.375
$ws.Range("J59").Value = 15793.375
$ws.Range("L59").Value = 15793.375
$ws.Range("N59").Value = -18083.375

$ws.Range("H60").Value = 7463.625
$ws.Range("J60").Value = 8201.286
$ws.Range("L60").Value = 8201.286
$ws.Range("N60").Value = -9223.286

$ws.Range("H61").Value = 9172.833000000001
$ws.Range("J61").Value = 9172.833000000001
$ws.Range("L61").Value = 9172.833000000001
$ws.Range("N61").Value = -9868.833000000001

$ws.Range("H68").Value = 18150.875
$ws.Range("J68").Value = 18150.875
$ws.Range("L68").Value = 18150.875
$ws.Range("N68").Value = -19648.875

$ws.Range("H71").Value = 18150.875
$ws.Range("J71").Value = 18150.875
$ws.Range("L71").Value = 54452.625
$ws.Range("N71").Value = -61940.625

$ws.Range("H74").Value = 18336.637
$ws.Range("J74").Value = 18336.637
$ws.Range("L74").Value = 18336.637
$ws.Range("N74").Value = -20084.637

$ws.Range("H77").Value = 18336.637
$ws.Range("J77").Value = 18336.637
$ws.Range("L77").Value = 55009.91099999999
$ws.Range("N77").Value = -63745.91099999999

$ws.Range("H132").Value = 2145.5293
$ws.Range("I132").Value = 1332.909
$ws.Range("K132").Value = 3998.727
$ws.Range("M132").Value = -1468.727

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1431724.4
$arr[0,1] = 2511.7058
$arr[0,2] = 3640507.5
$arr[0,3] = 7535.117400000001
$arr[0,4] = 10921522.5
$arr[0,5] = -5000.117400000001
$arr[0,6] = -10926592.5
$ws.Range("H134:N134").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 970153.5600000001
$arr[0,1] = 4085.3
$arr[0,2] = 2674979.8
$arr[0,3] = 12255.9
$arr[0,4] = 8024939.399999999
$arr[0,5] = -9705.900000000001
$arr[0,6] = -8030039.399999999
$ws.Range("H136:N136").Value = $arr

$ws = $wb.Worksheets.Item("CUL")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2404818.5
$arr[0,1] = 784.6667
$arr[0,2] = 5683046.5
$arr[0,3] = 2354.0001
$arr[0,4] = 17049139.5
$arr[0,5] = -2242.0001
$arr[0,6] = -17049363.5
$ws.Range("H5:N5").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 933.3838500000001
$arr[0,1] = 375
$arr[0,2] = 944.8969
$arr[0,3] = 1125
$arr[0,4] = 2834.6907
$arr[0,5] = 3915
$arr[0,6] = -12914.6907
$ws.Range("H131:N131").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2404818.5
$arr[0,1] = 784.6667
$arr[0,2] = 5683046.5
$arr[0,3] = 7062.0003
$arr[0,4] = 51147418.5
$arr[0,5] = -4527.0003
$arr[0,6] = -51152488.5
$ws.Range("H135:N135").Value = $arr

$ws = $wb.Worksheets.Item("GSM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 6666.375
$arr[0,1] = 11444.9
$arr[0,2] = 3253.1428
$arr[0,3] = 34334.7
$arr[0,4] = 9759.428400000001
$arr[0,5] = -31864.7
$arr[0,6] = -14699.4284
$ws.Range("H126:N126").Value = $arr

$ws.Range("H132").Value = 5884391.5
$ws.Range("I132").Value = 9525537
$ws.Range("K132").Value = 28576611
$ws.Range("M132").Value = -28574081
